$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.324.45'
$ws.Range('E2').Value = '  +0.16%  '

# Row 3
$ws.Range('D3').Value = '1.869.89'
$ws.Range('E3').Value = '  +0.26%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.14%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.08'
$ws.Range('E5').Value = '  -0.90%  '

# Row 6
$ws.Range('E6').Value = '  +0.13%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4704'
$ws.Range('E7').Value = '  +0.46%  '

# Row 8
$ws.Range('E8').Value = '  +0.33%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06567'
$ws.Range('E9').Value = '  +0.38%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.64'
$ws.Range('E10').Value = '  -3.17%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08025'
$ws.Range('E11').Value = '  +1.43%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '96.97'
$ws.Range('E12').Value = '  -0.78%  '

# Row 13
$ws.Range('D13').Value = '1.866.94'
$ws.Range('E13').Value = '  +0.09%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.111'
$ws.Range('E14').Value = '  -1.24%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6837'
$ws.Range('E15').Value = '  -0.08%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '269.38'
$ws.Range('E16').Value = '  -3.65%  '

# Row 17
$ws.Range('D17').Value = '30.288.94'
$ws.Range('E17').Value = '  +0.08%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.03'
$ws.Range('E18').Value = '  +2.58%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007608'
$ws.Range('E19').Value = '  +3.73%  '

# Row 20
$ws.Range('E20').Value = '  +0.12%  '

# Row 21
$ws.Range('D21').Value = '2.108.16'
$ws.Range('E21').Value = '  -0.19%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.14%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.280'
$ws.Range('E23').Value = '  -2.13%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.218'
$ws.Range('E24').Value = '  +0.70%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.429'
$ws.Range('E25').Value = '  +1.76%  '

# Row 26
$ws.Range('E26').Value = '  +0.39%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.90'
$ws.Range('E27').Value = '  -1.29%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.947'
$ws.Range('E28').Value = '  +0.53%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.371'
$ws.Range('E29').Value = '  -0.75%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09935'
$ws.Range('E30').Value = '  +1.21%  '

# Row 31
$ws.Range('E31').Value = '  -0.83%  '

# Row 32
$ws.Range('E32').Value = '  -1.20%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.075'
$ws.Range('E33').Value = '  +0.13%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04698'
$ws.Range('E34').Value = '  -1.10%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.139'
$ws.Range('E35').Value = '  +0.08%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7000'
$ws.Range('E36').Value = '  -1.61%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.700'
$ws.Range('E37').Value = '  -0.15%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01876'
$ws.Range('E38').Value = '  +0.06%  '

# Row 39
$ws.Range('E39').Value = '  +0.77%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.309'
$ws.Range('E40').Value = '  -0.06%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '71.83'
$ws.Range('E41').Value = '  -6.01%  '

# Row 42
$ws.Range('E42').Value = '  +0.36%  '

# Row 43
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4169'
$ws.Range('E43').Value = '  -0.52%  '

# Row 44
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8414'
$ws.Range('E44').Value = '  -1.11%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').Value = '  +0.13%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.88'
$ws.Range('E46').Value = '  -0.59%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.173'
$ws.Range('E47').Value = '  -1.57%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.050'
$ws.Range('E48').Value = '  -2.63%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '914.31'
$ws.Range('E49').Value = '  -5.58%  '

# Row 50
$ws.Range('E50').Value = '  +0.77%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05692'
$ws.Range('E51').Value = '  +0.83%  '
